# Flowchart-friendly rework of the "open_tasks" sheet (task 1b):
#  - explicitly note that "Double Check State Transition Diagram" depends on
#    "State Transition Diagram"
#  - widen the Task/Dependency columns so the longer text is readable
#  - tighten the row heights back to a single compact line
#  - clean up a couple of stray per-cell formats (B1 / A8) so they reuse the
#    same look as the rest of the header/task cells

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# 1. Fill in the now-explicit dependency for row 3.
$ws.Range("C3").Value = "State Transition Diagram"

# 2. Widen columns A and C so the (now longer) text fits without truncation.
$ws.Columns.Item(1).ColumnWidth = 49.17
$ws.Columns.Item(3).ColumnWidth = 49.17

# 3. Shrink the task rows back down to a single readable line each
#    (row 7 previously had to be very tall to show a wrapped note).
for ($r = 2; $r -le 7; $r++) {
    $ws.Rows.Item($r).RowHeight = 18.75
}

# 4. Normalize a couple of cells whose formatting had drifted (priority
#    header and the last task) back to the standard black-font look.
$ws.Range("B1").Font.Color = 0

$ws.Range("A8").HorizontalAlignment = -4131
$ws.Range("A8").Font.Color = 0
